$wb = $excel.ActiveWorkbook

# --- Sheet 1: hojaejemplo21 ---
$ws1 = $wb.Worksheets.Item("hojaejemplo21")
$ws1.Range("A1:A4").Copy($ws1.Range("B1:B4"))

# --- Sheet 2: hojaejemplo22 ---
$ws2 = $wb.Worksheets.Item("hojaejemplo22")
$ws2.Range("A1:A5").Copy($ws2.Range("B1:B5"))

# --- Sheet 3: hojaejemplo23 (data already present, only selection/view changes) ---
$ws3 = $wb.Worksheets.Item("hojaejemplo23")

# Update selections per sheet (without changing the active sheet yet)
$null = $ws1.Range("B1").Select()
$null = $ws3.Range("B1").Select()

# hojaejemplo22 becomes the active sheet/tab, with B2 selected
$null = $ws2.Activate()
$null = $ws2.Range("B2").Select()
